$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 29, whose C/D cells hold the
# "Total Duration:" summary label/value. A new clock-in/out entry needs to
# be appended as row 29, and the summary needs to move down to row 30 with
# an updated total.

# 1) Push the "Total Duration:" summary down from row 29 to the new row 30.
#    Copy first so the destination cells inherit the same cell style as the
#    source, then update the total value.
$ws.Cells.Item(29, 3).Copy($ws.Cells.Item(30, 3))
$ws.Cells.Item(29, 4).Copy($ws.Cells.Item(30, 4))
$ws.Cells.Item(30, 4).Value = "33.5 Hours"

# 2) Turn row 29 into the new timesheet entry. Prime the formatting of each
#    cell by copying down from row 28 (the last data row) so the new row
#    matches the look of the rest of the table, then fill in the values.
$ws.Cells.Item(28, 1).Copy($ws.Cells.Item(29, 1))
$ws.Cells.Item(28, 2).Copy($ws.Cells.Item(29, 2))
$ws.Cells.Item(28, 3).Copy($ws.Cells.Item(29, 3))
$ws.Cells.Item(28, 4).Copy($ws.Cells.Item(29, 4))

$ws.Cells.Item(29, 2).Value = "23:12:52"
$ws.Cells.Item(29, 3).Value = "23:55:46"
$ws.Cells.Item(29, 4).Value = "0.71 Hours"

# A29 needs special handling: assigning a date-shaped string like
# "2026-02-18" straight to .Value gets auto-detected and converted into a
# date serial number (and resets the cell's style in the process).
# Entering it as a "=<quoted text>" formula avoids that auto-detection
# (formula entry isn't re-interpreted as a date), and then a values-only
# paste bakes the formula's result back down into a plain literal text
# value while keeping the cell's style untouched.
$ws.Cells.Item(29, 1).Value = '="2026-02-18"'
$ws.Cells.Item(29, 1).Copy()
$ws.Cells.Item(29, 1).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = $false
